$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: run a simple Find & Replace over the whole document body (story).
# wdReplaceAll semantics aren't needed here since every target string is
# unique; MatchCase=$true, Wrap=wdFindContinue(1), Replace=wdReplaceOne(2)
# ---------------------------------------------------------------------------
function Replace-Text($findText, $replaceText) {
    $rng = $d.Content
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# 1) "タイムライン、リソースの割り当て、コンティンジェンシー 計画を含めます。"
Replace-Text "タイムライン、リソースの割り当て、コンティンジェンシー 計画を含めます。" "タイムライン、リソース割り当て、コンティンジェンシー計画を含めます。"

# 2) training-program sentence rewording
Replace-Text "Contoso CipherGuard Sentinel X7 の使用方法と保守方法について、すべてのユーザーと管理者が適切にトレーニングされるように、トレーニング プログラムを開発して実装します。" "Contoso CipherGuard Sentinel X7 を使用して維持する方法について、すべてのユーザーと管理者が適切にトレーニングされるように、トレーニング プログラムを開発して実装します。"

# 3) communication-plan sentence: 関係者 -> 利害関係者
Replace-Text "Contoso CipherGuard Sentinel X7 のデプロイについてすべての関係者に確実に通知されるように、コミュニケーション計画を策定して実装します。" "Contoso CipherGuard Sentinel X7 のデプロイについてすべての利害関係者に確実に通知されるように、コミュニケーション計画を策定して実装します。"

# 4) "文書化とレポート作成" is followed by a separate, non-bold ": " run;
#    trim the trailing space from just that run without merging it with the
#    preceding bold run.
$leadRun = $d.Content
$leadRun.Find.Execute("文書化とレポート作成", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$colonRun = $d.Range($leadRun.End, $leadRun.End + 2)
$colonRun.Text = ":"

# 5) documentation/report sentence: 展開 -> デプロイ
Replace-Text "Contoso CipherGuard Sentinel X7 の展開に関するすべての関連情報が適切に文書化および報告されるように、ドキュメントとレポート計画を作成して実装します。" "Contoso CipherGuard Sentinel X7 のデプロイに関するすべての関連情報が適切に文書化および報告されるように、ドキュメントとレポート計画を作成して実装します。"

# 6) "プロジェクトのタイムライン。" heading run:
#    - becomes bold
#    - loses the trailing "。"
#    - gets a new, non-bold ":" run appended (cloning formatting from the
#      "文書化とレポート作成:" colon run fixed up in step 4 above, which
#      already carries the exact w:val="0" boolean property set Word uses)
$timelineRun = $d.Content
$timelineRun.Find.Execute("プロジェクトのタイムライン。", $true, $false, $false, $false, $false, $true, 1, $false, "プロジェクトのタイムライン", 2) | Out-Null
$timelineRun.Font.Bold = 1
$afterTimeline = $timelineRun.End

$tplLead = $d.Content
$tplLead.Find.Execute("文書化とレポート作成", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tplColon = $d.Range($tplLead.End, $tplLead.End + 1)

$insertionPoint = $d.Range($afterTimeline, $afterTimeline)
$insertionPoint.FormattedText = $tplColon.FormattedText

# 7) "プロジェクトタイムライン" -> "プロジェクト タイムライン" (add a space)
Replace-Text "Contoso CipherGuard Sentinel X7 のデプロイに関連するすべてのタスクとアクティビティが期限内に完了するように、プロジェクトタイムラインを開発して実装します。" "Contoso CipherGuard Sentinel X7 のデプロイに関連するすべてのタスクとアクティビティが期限内に完了するように、プロジェクト タイムラインを開発して実装します。"

# 8) "これは単なるサンプル 計画であり、..." -> remove the space before 計画.
#    This run sits right before another run ("新しいネットワーク...") that
#    happens to share identical run formatting; a plain Find/Replace on run 1
#    causes the engine to silently coalesce it with run 2. Re-split them
#    afterwards by cloning run 2's own FormattedText back onto itself, which
#    forces a fresh run boundary while reproducing the exact explicit
#    w:val="0" boolean property set Word already uses on both runs.
Replace-Text "これは単なるサンプル 計画であり、組織の特定のニーズと要件を満たすために調整する必要がある場合があります。" "これは単なるサンプル計画であり、組織の特定のニーズと要件を満たすために調整する必要がある場合があります。"

$followingRun = $d.Content
$followingRun.Find.Execute("新しいネットワーク セキュリティ製品をデプロイするときは、常に業界の専門家と相談し、確立されたベスト プラクティスに従うことをお勧めします。", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$followingRun.FormattedText = $followingRun.FormattedText

# 9) Header text: "目を向ける際、" -> "AI を活用"
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdrRange = $hdr.Range
$hdrRange.Find.Execute("目を向ける際、", $true, $false, $false, $false, $false, $true, 1, $false, "AI を活用", 2) | Out-Null

Write-Output "edits applied"
